$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the data range so numeric-looking strings
# (e.g. "1.006") are not reinterpreted as numbers, matching the
# original inlineStr/shared-string text cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.409.30"
$ws.Range("E2").Value = "  +8.77%  "

$ws.Range("D3").Value = "1.677.77"
$ws.Range("E3").Value = "  +4.50%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "305.99"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "0.3696"
$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").Value = "0.3430"
$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("D9").Value = "47.55"
$ws.Range("E9").Value = "  +12.34%  "

$ws.Range("D10").Value = "1.162"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").Value = "0.07229"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").Value = "6.100"
$ws.Range("E13").Value = "  +3.16%  "

$ws.Range("D14").Value = "20.11"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "6.736"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").Value = "1.675.34"
$ws.Range("E16").Value = "  +4.29%  "

$ws.Range("D17").Value = "0.00001099"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").Value = "0.9992"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "0.06661"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "80.68"
$ws.Range("E20").Value = "  +3.61%  "

$ws.Range("D21").Value = "16.40"
$ws.Range("E21").Value = "  +2.35%  "

$ws.Range("D22").Value = "6.094"
$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").Value = "12.11"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").Value = "24.361.77"
$ws.Range("E24").Value = "  +8.28%  "

$ws.Range("D25").Value = "2.431"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.640"
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "152.16"
$ws.Range("E27").Value = "  +1.76%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.41"
$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "1.860.93"
$ws.Range("E29").Value = "  +4.06%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "127.74"
$ws.Range("E30").Value = "  +4.14%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "6.269"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "4.058"
$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.9682"
$ws.Range("E33").Value = "  +1.99%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.08421"
$ws.Range("E34").Value = "  +2.24%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.680"
$ws.Range("E35").Value = "  -0.91%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "12.20"
$ws.Range("E36").Value = "  +1.38%  "

$ws.Range("D37").Value = "0.06398"
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "5.307"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.734"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02312"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.232"
$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.2088"
$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.6092"
$ws.Range("E43").Value = "  +3.19%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.754"
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.01"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5876"
$ws.Range("E47").Value = "  +3.33%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "126.63"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.014"
$ws.Range("E49").Value = "  +2.63%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07149"
$ws.Range("E50").Value = "  +4.88%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "75.62"
$ws.Range("E51").Value = "  +2.87%  "

# Restore the default (Normal) style on the data range so no stray
# text-format styling is left behind on the cells.
$dataRange.Style = "Normal"
